# Update the EPEX Spot prices workbook with the latest day of data.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Prix Spot": append a new "21-jun" column (H) with hourly prices
# ------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("G1").Copy()
$wsPrix.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$wsPrix.Range("H1").Value = "21-jun"

$prixSpotValues = @(
    118.53,
    110.48,
    108.15,
    98.03,
    81.62,
    92.27,
    98.79000000000001,
    90.73999999999999,
    75.93000000000001,
    23.4,
    2.5,
    0,
    0,
    0,
    0,
    0,
    0.1,
    60.4,
    98.64,
    127.53,
    141.99,
    134.93,
    138.81,
    119.48
)

for ($i = 0; $i -lt $prixSpotValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 8).Value = $prixSpotValues[$i]
}

# ------------------------------------------------------------------
# Sheet "Gaz": append the new day's last price
# ------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force the date text to stay a plain text string (not auto-converted to
# a date serial number) while not leaving a lingering date number format.
$wsGaz.Range("A5").NumberFormat = "@"
$wsGaz.Range("A5").Value = "2025-06-19"
$wsGaz.Range("A5").ClearFormats()
$wsGaz.Range("B5").Value = 40.425

# ------------------------------------------------------------------
# Sheet "CO2": append the new day's last price
# ------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A5").NumberFormat = "@"
$wsCo2.Range("A5").Value = "2025-06-19"
$wsCo2.Range("A5").ClearFormats()
$wsCo2.Range("B5").Value = 71
